# Update MSME Country Indicators - Germany Summary values
# Source Type: SME Associations (Most Widely Used) section (rows 32, 34, 36)
#
# These figures are stored as plain text (General format, shared-string
# cells) rather than numbers, so assigning a numeric-looking string
# directly would make Excel auto-convert the cell to a Number. Temporarily
# switch the cell to Text format, assign the new value so it is kept as a
# string, then restore the cell style so formatting is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Employment (% of total) row: Micro=18.5, SMEs, MSMEs
Set-TextValue "C32" "43.74"
Set-TextValue "D32" "62.25"

# Enterprises density (per 1000 people) row: Micro, SMEs=4.7, MSMEs
Set-TextValue "B34" "21.56"
Set-TextValue "D34" "26.25"

# Enterprises (% of total) row: Micro, SMEs=17.8, MSMEs=99.5
Set-TextValue "B36" "81.71"
